$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, shifting existing rows 22+ down by one.
# (formatting is inherited from the row above, matching neighboring data rows)
$ws.Rows.Item(22).Insert()

# Populate the new row with the "description is too long" translation triple,
# in B, C, A order so new shared strings land in the same order as the target file.
$ws.Range("B22").Value = "your description is so long!"
$ws.Range("C22").Value = "توضیحات شما بسیار طولانی است!"
$ws.Range("A22").Value = "your_description_is_so_long"

# Update the view state to match: scroll back to top, select the new row's first cell.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("A22").Select()
